$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.325956666666666
$ws.Range("H2").Value = 18.97787
$ws.Range("I2").Value = 0.4468357575736242
$ws.Range("J2").Value = 0.4592138460625664
$ws.Range("M2").Value = 2.535712666666667
$ws.Range("N2").Value = 7.607138
$ws.Range("O2").Value = 0.04494879354621957
$ws.Range("P2").Value = 0.05070282964779482
$ws.Range("Q2").Value = 16.04080844845111
$ws.Range("R2").Value = 144.36727603606
$ws.Range("S2").Value = 0.02008472821624545
$ws.Range("T2").Value = 0.02328344140881897

$ws.Range("G3").Value = 6.325956666666666
$ws.Range("H3").Value = 18.97787
$ws.Range("I3").Value = 0.4468357575736242
$ws.Range("J3").Value = 0.4592138460625664
$ws.Range("O3").Value = 0.5715421877013505
$ws.Range("P3").Value = 0.6447070965264385
$ws.Range("Q3").Value = 203.9654021792344
$ws.Range("R3").Value = 1835.68861961311
$ws.Range("S3").Value = 0.2553854864268195
$ws.Range("T3").Value = 0.2960584253797361

$ws.Range("G4").Value = 6.325956666666666
$ws.Range("H4").Value = 18.97787
$ws.Range("I4").Value = 0.4468357575736242
$ws.Range("J4").Value = 0.4592138460625664
$ws.Range("M4").Value = 1.538811333333333
$ws.Range("N4").Value = 4.616434
$ws.Range("O4").Value = 0.02727742533206951
$ws.Range("P4").Value = 0.03076929413956839
$ws.Range("Q4").Value = 9.734453812842222
$ws.Range("R4").Value = 87.61008431558
$ws.Range("S4").Value = 0.01218852901291325
$ws.Range("T4").Value = 0.01412968590246158

$ws.Range("G5").Value = 6.325956666666666
$ws.Range("H5").Value = 18.97787
$ws.Range("I5").Value = 0.4468357575736242
$ws.Range("J5").Value = 0.4592138460625664
$ws.Range("M5").Value = 19.206297
$ws.Range("N5").Value = 38.412594
$ws.Range("O5").Value = 0.3404565075487166
$ws.Range("P5").Value = 0.2560262755732715
$ws.Range("Q5").Value = 121.49820254913
$ws.Range("R5").Value = 728.9892152947799
$ws.Range("S5").Value = 0.1521281414714011
$ws.Range("T5").Value = 0.1175708106990765

$ws.Range("G6").Value = 6.325956666666666
$ws.Range("H6").Value = 18.97787
$ws.Range("I6").Value = 0.4468357575736242
$ws.Range("J6").Value = 0.4592138460625664
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8899256666666666
$ws.Range("N6").Value = 2.669777
$ws.Range("O6").Value = 0.0157750858716439
$ws.Range("P6").Value = 0.01779450411292666
$ws.Range("Q6").Value = 5.629631203887778
$ws.Range("R6").Value = 50.66668083499
$ws.Range("S6").Value = 0.007048872446244978
$ws.Range("T6").Value = 0.008171482672473208

$ws.Range("G7").Value = 6.686451000000001
$ws.Range("I7").Value = 0.4722993778644153
$ws.Range("J7").Value = 0.4853828506917099
$ws.Range("M7").Value = 2.535712666666667
$ws.Range("N7").Value = 7.607138
$ws.Range("O7").Value = 0.04494879354621957
$ws.Range("P7").Value = 0.05070282964779482
$ws.Range("Q7").Value = 16.954918495746
$ws.Range("R7").Value = 152.594266461714
$ws.Range("S7").Value = 0.02122928722763555
$ws.Range("T7").Value = 0.02461028399258279

$ws.Range("G8").Value = 6.686451000000001
$ws.Range("I8").Value = 0.4722993778644153
$ws.Range("J8").Value = 0.4853828506917099
$ws.Range("O8").Value = 0.5715421877013505
$ws.Range("P8").Value = 0.6447070965264385
$ws.Range("S8").Value = 0.2699390196746147
$ws.Range("T8").Value = 0.3129297683731781

$ws.Range("G9").Value = 6.686451000000001
$ws.Range("I9").Value = 0.4722993778644153
$ws.Range("J9").Value = 0.4853828506917099
$ws.Range("M9").Value = 1.538811333333333
$ws.Range("N9").Value = 4.616434
$ws.Range("O9").Value = 0.02727742533206951
$ws.Range("P9").Value = 0.03076929413956839
$ws.Range("Q9").Value = 10.289186578578
$ws.Range("R9").Value = 92.602679207202
$ws.Range("S9").Value = 0.01288311101407947
$ws.Range("T9").Value = 0.01493488770323543

$ws.Range("G10").Value = 6.686451000000001
$ws.Range("I10").Value = 0.4722993778644153
$ws.Range("J10").Value = 0.4853828506917099
$ws.Range("M10").Value = 19.206297
$ws.Range("N10").Value = 38.412594
$ws.Range("O10").Value = 0.3404565075487166
$ws.Range("P10").Value = 0.2560262755732715
$ws.Range("Q10").Value = 128.421963781947
$ws.Range("R10").Value = 770.531782691682
$ws.Range("S10").Value = 0.1607973967051505
$ws.Range("T10").Value = 0.1242707634897358

$ws.Range("G11").Value = 6.686451000000001
$ws.Range("I11").Value = 0.4722993778644153
$ws.Range("J11").Value = 0.4853828506917099
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.8899256666666666
$ws.Range("N11").Value = 2.669777
$ws.Range("O11").Value = 0.0157750858716439
$ws.Range("P11").Value = 0.01779450411292666
$ws.Range("Q11").Value = 5.950444363809001
$ws.Range("R11").Value = 53.553999274281
$ws.Range("S11").Value = 0.007450563242935143
$ws.Range("T11").Value = 0.008637147132977698

$ws.Range("G12").Value = 1.1448225
$ws.Range("H12").Value = 2.289645
$ws.Range("I12").Value = 0.08086486456196039
$ws.Range("J12").Value = 0.05540330324572383
$ws.Range("M12").Value = 2.535712666666667
$ws.Range("N12").Value = 7.607138
$ws.Range("O12").Value = 0.04494879354621957
$ws.Range("P12").Value = 0.05070282964779482
$ws.Range("Q12").Value = 2.902940914335
$ws.Range("R12").Value = 17.41764548601
$ws.Range("S12").Value = 0.003634778102338565
$ws.Range("T12").Value = 0.002809104246393053

$ws.Range("G13").Value = 1.1448225
$ws.Range("H13").Value = 2.289645
$ws.Range("I13").Value = 0.08086486456196039
$ws.Range("J13").Value = 0.05540330324572383
$ws.Range("O13").Value = 0.5715421877013505
$ws.Range("P13").Value = 0.6447070965264385
$ws.Range("Q13").Value = 36.9120741636975
$ws.Range("R13").Value = 221.472444982185
$ws.Range("S13").Value = 0.04621768159991625
$ws.Range("T13").Value = 0.03571890277352442

$ws.Range("G14").Value = 1.1448225
$ws.Range("H14").Value = 2.289645
$ws.Range("I14").Value = 0.08086486456196039
$ws.Range("J14").Value = 0.05540330324572383
$ws.Range("M14").Value = 1.538811333333333
$ws.Range("N14").Value = 4.616434
$ws.Range("O14").Value = 0.02727742533206951
$ws.Range("P14").Value = 0.03076929413956839
$ws.Range("Q14").Value = 1.761665837655
$ws.Range("R14").Value = 10.56999502593
$ws.Range("S14").Value = 0.002205785305076789
$ws.Range("T14").Value = 0.001704720533871381

$ws.Range("G15").Value = 1.1448225
$ws.Range("H15").Value = 2.289645
$ws.Range("I15").Value = 0.08086486456196039
$ws.Range("J15").Value = 0.05540330324572383
$ws.Range("M15").Value = 19.206297
$ws.Range("N15").Value = 38.412594
$ws.Range("O15").Value = 0.3404565075487166
$ws.Range("P15").Value = 0.2560262755732715
$ws.Range("Q15").Value = 21.9878009472825
$ws.Range("R15").Value = 87.95120378913001
$ws.Range("S15").Value = 0.02753096937216501
$ws.Range("T15").Value = 0.01418470138445921

$ws.Range("G16").Value = 1.1448225
$ws.Range("H16").Value = 2.289645
$ws.Range("I16").Value = 0.08086486456196039
$ws.Range("J16").Value = 0.05540330324572383
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.8899256666666666
$ws.Range("N16").Value = 2.669777
$ws.Range("O16").Value = 0.0157750858716439
$ws.Range("P16").Value = 0.01779450411292666
$ws.Range("Q16").Value = 1.0188069265275
$ws.Range("R16").Value = 6.112841559165
$ws.Range("S16").Value = 0.001275650182463779
$ws.Range("T16").Value = 0.0009858743074757557

Write-Output "applied updates"